$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '88.434.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.18%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.276.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.79%  '

# Row 4
$ws.Range("E4").Value = '  +0.21%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.98%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '630.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.77%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.392'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +22.83%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.690'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +16.68%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.12%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.272.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.87%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.579'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.51%  '

# Row 12
$ws.Range("E12").Value = '  +10.81%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000264'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.25%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.98%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.880.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.52%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.81%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.157.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.12%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.263.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.87%  '

# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.65%  '

# Row 20
$ws.Range("B20").Value = 'SuiNetwork'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.52%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.93'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.33%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.23%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.88%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.66%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.421.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.14%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '77.02'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.09%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000134'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.59%  '

# Row 30
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.183'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.55%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.08%  '

# Row 33
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '571.74'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.30%  '

# Row 34
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '8.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.05%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.41'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.74%  '

# Row 36
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.97'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.49%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.15'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.91%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.139'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.93%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.68'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.00%  '

# Row 40
$ws.Range("E40").Value = '  +0.20%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '21.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.86%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.398'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.58%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.52%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.53%  '

# Row 45
$ws.Range("E45").Value = '  +0.13%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '152.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.52%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '180.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.22%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.12%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.30'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.38%  '

# Row 50
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.25'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.33%  '

# Row 51
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.125'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.57%  '
